$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the two "not tracked" cells entirely (content + formatting),
# so the shared string "not tracked" becomes unused and is dropped on save.
$ws.Range("C4").Clear()
$ws.Range("C5").Clear()

# New note next to "Doc" row explaining same-as-click tracking behaviour.
$ws.Range("E2").Value = "same as click"

# New row documenting that in-page jumping isn't tracked.
$ws.Range("B9").Value = "in-page jumping is not tracked"

# Column G widened (best-fit) to accommodate the longer tracked text.
$ws.Columns.Item(7).ColumnWidth = 27

# Selection / view state as recorded by the author when saving.
$ws.Range("F6").Select()
